# Auto-generated Excel COM-interop script applying the Chocobo_Profits.xlsx diff
# Source: scheduled runner data refresh (currentAveragePrice / Leve profit columns)
$wb = $excel.ActiveWorkbook

# --- ALC row 129 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1001.0455
$ws.Range("J129").Value = 1012.72095
$ws.Range("L129").Value = 3038.16285
$ws.Range("N129").Value = -13038.16285

# --- ALC row 137 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 3754.3667
$ws.Range("I137").Value = 2380.2856
$ws.Range("J137").Value = 6960.5557
$ws.Range("K137").Value = 7140.8568
$ws.Range("L137").Value = 20881.6671
$ws.Range("M137").Value = -4590.8568
$ws.Range("N137").Value = -25981.6671

# --- ARM row 132 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 2425.587
$ws.Range("I132").Value = 1273.5217
$ws.Range("J132").Value = 3577.652
$ws.Range("K132").Value = 3820.5651
$ws.Range("L132").Value = 10732.956
$ws.Range("M132").Value = -1290.5651
$ws.Range("N132").Value = -15792.956

# --- BSM row 19 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 39333
$ws.Range("J19").Value = 39333
$ws.Range("L19").Value = 39333
$ws.Range("N19").Value = -39679

# --- BSM row 86 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2236.7646
$ws.Range("I86").Value = 1877.0834
$ws.Range("J86").Value = 3100
$ws.Range("K86").Value = 1877.0834
$ws.Range("L86").Value = 3100
$ws.Range("M86").Value = -754.0834
$ws.Range("N86").Value = -5346

# --- BSM row 89 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 2236.7646
$ws.Range("I89").Value = 1877.0834
$ws.Range("J89").Value = 3100
$ws.Range("K89").Value = 9385.416999999999
$ws.Range("L89").Value = 15500
$ws.Range("M89").Value = -3769.416999999999
$ws.Range("N89").Value = -26732

# --- BSM row 134 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3547.4473
$ws.Range("I134").Value = 2079.875
$ws.Range("J134").Value = 6063.2856
$ws.Range("K134").Value = 6239.625
$ws.Range("L134").Value = 18189.8568
$ws.Range("M134").Value = -3704.625
$ws.Range("N134").Value = -23259.8568

# --- CRP row 31 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3370.5264
$ws.Range("I31").Value = 1531.2
$ws.Range("J31").Value = 6907.6924
$ws.Range("K31").Value = 1531.2
$ws.Range("L31").Value = 6907.6924
$ws.Range("M31").Value = -1236.2
$ws.Range("N31").Value = -7497.6924

# --- CRP row 34 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 3370.5264
$ws.Range("I34").Value = 1531.2
$ws.Range("J34").Value = 6907.6924
$ws.Range("K34").Value = 1531.2
$ws.Range("L34").Value = 6907.6924
$ws.Range("M34").Value = -1329.2
$ws.Range("N34").Value = -7311.6924

# --- CRP row 68 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 47474.332
$ws.Range("J68").Value = 47474.332
$ws.Range("L68").Value = 47474.332
$ws.Range("N68").Value = -48972.332

# --- CRP row 71 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 47474.332
$ws.Range("J71").Value = 47474.332
$ws.Range("L71").Value = 142422.996
$ws.Range("N71").Value = -149910.996

# --- CUL row 81 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 1504.3334
$ws.Range("I81").Value = 1006.5
$ws.Range("J81").Value = 2500
$ws.Range("K81").Value = 3019.5
$ws.Range("L81").Value = 7500
$ws.Range("M81").Value = -1896.5
$ws.Range("N81").Value = -9746

# --- CUL row 84 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H84").Value = 1504.3334
$ws.Range("I84").Value = 1006.5
$ws.Range("J84").Value = 2500
$ws.Range("K84").Value = 9058.5
$ws.Range("L84").Value = 22500
$ws.Range("M84").Value = -3442.5
$ws.Range("N84").Value = -33732

# --- CUL row 124 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H124").Value = 3454.5454
$ws.Range("J124").Value = 3454.5454
$ws.Range("L124").Value = 10363.6362
$ws.Range("N124").Value = -20183.6362

# --- CUL row 125 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H125").Value = 1000
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# --- CUL row 129 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 2629.5264
$ws.Range("I129").Value = 0
$ws.Range("J129").Value = 2629.5264
$ws.Range("K129").Value = 0
$ws.Range("L129").Value = 7888.5792
$ws.Range("M129").ClearContents()
$ws.Range("N129").Value = -17888.5792

# --- CUL row 130 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H130").Value = 2846.6667
$ws.Range("I130").Value = 1926.6666
$ws.Range("J130").Value = 3766.6667
$ws.Range("K130").Value = 5779.9998
$ws.Range("L130").Value = 11300.0001
$ws.Range("M130").Value = -759.9997999999996
$ws.Range("N130").Value = -21340.0001

# --- CUL row 131 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 7579756
$ws.Range("I131").Value = 26327198
$ws.Range("J131").Value = 1002.8511
$ws.Range("K131").Value = 78981594
$ws.Range("L131").Value = 3008.5533
$ws.Range("M131").Value = -78976554
$ws.Range("N131").Value = -13088.5533

# --- CUL row 132 ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 3256.5557
$ws.Range("J132").Value = 4068.2122
$ws.Range("L132").Value = 36613.9098
$ws.Range("N132").Value = -41673.9098

# --- GSM row 98 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 15821.5
$ws.Range("J98").Value = 15821.5
$ws.Range("L98").Value = 15821.5
$ws.Range("N98").Value = -21811.5

# --- GSM row 124 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H124").Value = 43113.332
$ws.Range("J124").Value = 43113.332
$ws.Range("L124").Value = 43113.332
$ws.Range("N124").Value = -52933.332

# --- GSM row 126 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 3940.649
$ws.Range("I126").Value = 2843.5686
$ws.Range("J126").Value = 5241.8374
$ws.Range("K126").Value = 8530.7058
$ws.Range("L126").Value = 15725.5122
$ws.Range("M126").Value = -6060.7058
$ws.Range("N126").Value = -20665.5122

# --- LTW row 7 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5899.5
$ws.Range("I7").Value = 3471.2856
$ws.Range("J7").Value = 8327.714
$ws.Range("K7").Value = 3471.2856
$ws.Range("L7").Value = 8327.714
$ws.Range("M7").Value = -3359.2856
$ws.Range("N7").Value = -8551.714

# --- LTW row 62 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()

# --- LTW row 65 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()

# --- LTW row 126 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5899.5
$ws.Range("I126").Value = 3471.2856
$ws.Range("J126").Value = 8327.714
$ws.Range("K126").Value = 10413.8568
$ws.Range("L126").Value = 24983.142
$ws.Range("M126").Value = -7943.856800000001
$ws.Range("N126").Value = -29923.142

# --- LTW row 140 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 67000
$ws.Range("J140").Value = 67000
$ws.Range("L140").Value = 67000
$ws.Range("N140").Value = -77360

# --- WVR row 70 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 23000
$ws.Range("J70").Value = 23000
$ws.Range("L70").Value = 23000
$ws.Range("N70").Value = -23630

# --- WVR row 73 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 23000
$ws.Range("J73").Value = 23000
$ws.Range("L73").Value = 23000
$ws.Range("N73").Value = -25184

